$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(106).Insert()

$ws.Range("A106").Value = "T12"
$ws.Range("B106").Value = "early_icu"
$ws.Range("C106").Value = "Time measurements"
$ws.Range("D106").Value = "Early ICU admission (within 48 hours)"

$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:E113"))

$excel.ActiveWindow.ScrollRow = 93
$ws.Range("E106").Select()
